$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 33: add the commit message text and hours value
$ws.Range("C33").Value = "game source update & potions added"
$ws.Range("G33").Value = 2

# Extend the Total(h) SUM formula to include the new row 33
$ws.Range("G39").Formula = "=SUM(G4:G33)"

# Update the view: scroll position and selection, matching the author's final state
$win = $excel.ActiveWindow
$win.ScrollRow = 28
$win.ScrollColumn = 1
$ws.Range("D34").Select()
